# Sync SharePoint -> GitHub (Liste Agents.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update Zone values for rows 2-4 (PRM4/PRM11/PRM3 -> PRM2)
$ws.Range("B2").Value = "PRM2"
$ws.Range("B3").Value = "PRM2"
$ws.Range("B4").Value = "PRM2"

# Row 3 "EPI" column corrected from NON to OUI
$ws.Range("F3").Value = "OUI"

# Update active selection to F3
$ws.Range("F3").Select()
